$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 73 ---
$ws.Range("C73").Value = "20:00:04"
$ws.Range("D73").Value = 1.13
$ws.Range("F73").Value = 69
$ws.Range("H73").Value = 2.24

# --- Update existing row 74 ---
$ws.Range("C74").Value = "20:30:05"
$ws.Range("D74").Value = 0.8100000000000001
$ws.Range("F74").Value = 72
$ws.Range("H74").Value = 1.34

# --- Append new row 75 ---
$ws.Range("A75").Value = "December"
$ws.Range("B75").Value = 17
$ws.Range("C75").Value = "21:00:05"
$ws.Range("D75").Value = 0.8100000000000001
$ws.Range("E75").Value = 0
$ws.Range("F75").Value = 76
$ws.Range("G75").Value = 1031
$ws.Range("H75").Value = 0.89

# --- Append new row 76 ---
$ws.Range("A76").Value = "December"
$ws.Range("B76").Value = 17
$ws.Range("C76").Value = "21:30:07"
$ws.Range("D76").Value = 0.8100000000000001
$ws.Range("E76").Value = 0
$ws.Range("F76").Value = 77
$ws.Range("G76").Value = 1031
$ws.Range("H76").Value = 0.45

# --- Append new row 77 ---
$ws.Range("A77").Value = "December"
$ws.Range("B77").Value = 17
$ws.Range("C77").Value = "22:00:05"
$ws.Range("D77").Value = 0.8100000000000001
$ws.Range("E77").Value = 0
$ws.Range("F77").Value = 78
$ws.Range("G77").Value = 1031
$ws.Range("H77").Value = 0.89

# --- Append new row 78 ---
$ws.Range("A78").Value = "December"
$ws.Range("B78").Value = 17
$ws.Range("C78").Value = "22:30:04"
$ws.Range("D78").Value = 1.13
$ws.Range("E78").Value = 0
$ws.Range("F78").Value = 80
$ws.Range("G78").Value = 1031
$ws.Range("H78").Value = 2.05

# --- Append new row 79 ---
$ws.Range("A79").Value = "December"
$ws.Range("B79").Value = 17
$ws.Range("C79").Value = "23:00:05"
$ws.Range("D79").Value = 1.5
$ws.Range("E79").Value = 0
$ws.Range("F79").Value = 81
$ws.Range("G79").Value = 1031
$ws.Range("H79").Value = 2.05
